# Apply weekly crypto price/volume refresh (GitHub Actions data pull).
# Diff-derived cell overwrites on Sheet1: columns D (Price) and E (Volume 1h),
# plus a B/C/D/E swap for rows 44-45 (ApeXProtocol <-> ThetaToken).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '70.246.72'
# Row 3: Ethereum
$ws.Range('D3').Value = '3.585.10'
$ws.Range('E3').Value = '  -1.32%  '
# Row 4: TetherUSD
$ws.Range('E4').Value = '  -0.02%  '
# Row 5: BNB
$ws.Range('D5').Value = '''578.85'
$ws.Range('E5').Value = '  -2.21%  '
# Row 6: Solana
$ws.Range('D6').Value = '''188.06'
$ws.Range('E6').Value = '  -3.58%  '
# Row 7: LidoStakedEther
$ws.Range('D7').Value = '3.580.08'
$ws.Range('E7').Value = '  -1.28%  '
# Row 8: XRP
$ws.Range('D8').Value = '''0.622'
$ws.Range('E8').Value = '  -3.49%  '
# Row 9: USDC
$ws.Range('E9').Value = '  +0.02%  '
# Row 10: Dogecoin
$ws.Range('D10').Value = '''0.186'
$ws.Range('E10').Value = '  +1.15%  '
# Row 11: Cardano
$ws.Range('D11').Value = '''0.656'
$ws.Range('E11').Value = '  -3.26%  '
# Row 12: Avalanche
$ws.Range('D12').Value = '''54.94'
$ws.Range('E12').Value = '  -5.26%  '
# Row 13: ShibaInu
$ws.Range('D13').Value = '''0.0000307'
$ws.Range('E13').Value = '  -1.40%  '
# Row 14: Polkadot
$ws.Range('D14').Value = '''9.60'
$ws.Range('E14').Value = '  -3.41%  '
# Row 15: WrappedliquidstakedEther2.0
$ws.Range('D15').Value = '4.160.59'
$ws.Range('E15').Value = '  -1.32%  '
# Row 16: Chainlink
$ws.Range('E16').Value = '  -3.56%  '
# Row 17: WrappedEther
$ws.Range('D17').Value = '3.586.83'
$ws.Range('E17').Value = '  -1.33%  '
# Row 18: WrappedBTC
$ws.Range('D18').Value = '70.148.23'
$ws.Range('E18').Value = '  -0.84%  '
# Row 19: Uniswap
$ws.Range('D19').Value = '''12.55'
$ws.Range('E19').Value = '  -1.52%  '
# Row 20: TRON
$ws.Range('E20').Value = '  -1.28%  '
# Row 21: Polygon
$ws.Range('E21').Value = '  -2.21%  '
# Row 22: BitcoinCash
$ws.Range('D22').Value = '''493.07'
$ws.Range('E22').Value = '  +0.99%  '
# Row 23: InternetComputer(DFINITY)
$ws.Range('D23').Value = '''19.72'
$ws.Range('E23').Value = '  +1.16%  '
# Row 24: Toncoin
$ws.Range('D24').Value = '''4.98'
$ws.Range('E24').Value = '  -4.40%  '
# Row 25: Litecoin
$ws.Range('D25').Value = '''97.64'
$ws.Range('E25').Value = '  +6.97%  '
# Row 26: PancakeSwap
$ws.Range('E26').Value = '  -1.81%  '
# Row 27: RenderToken
$ws.Range('D27').Value = '''11.45'
$ws.Range('E27').Value = '  +0.25%  '
# Row 28: ImmutableX
$ws.Range('E28').Value = '  -5.81%  '
# Row 29: Filecoin
$ws.Range('D29').Value = '''9.38'
$ws.Range('E29').Value = '  -1.84%  '
# Row 30: NEARProtocol
$ws.Range('D30').Value = '''7.72'
$ws.Range('E30').Value = '  -3.24%  '
# Row 31: EthereumClassic
$ws.Range('D31').Value = '''31.78'
$ws.Range('E31').Value = '  -3.15%  '
# Row 32: Cosmos
$ws.Range('D32').Value = '''12.28'
$ws.Range('E32').Value = '  -0.04%  '
# Row 33: OKB
$ws.Range('D33').Value = '''65.74'
$ws.Range('E33').Value = '  -0.69%  '
# Row 34: Hedera
$ws.Range('E34').Value = '  -5.19%  '
# Row 35: Bittensor
$ws.Range('D35').Value = '''575.72'
$ws.Range('E35').Value = '  -5.84%  '
# Row 36: Fetch.AI
$ws.Range('D36').Value = '''3.19'
$ws.Range('E36').Value = '  +11.11%  '
# Row 37: InjectiveProtocol
$ws.Range('D37').Value = '''38.93'
$ws.Range('E37').Value = '  -3.95%  '
# Row 38: TheGraph
$ws.Range('E38').Value = '  -1.36%  '
# Row 39: Dai
$ws.Range('E39').Value = '  -0.02%  '
# Row 40: PEPE
$ws.Range('E40').Value = '  -5.17%  '
# Row 41: Stacks
$ws.Range('E41').Value = '  -2.60%  '
# Row 42: dogwifhat
$ws.Range('D42').Value = '''3.19'
$ws.Range('E42').Value = '  +1.54%  '
# Row 43: Kaspa
$ws.Range('E43').Value = '  -8.24%  '
# Row 44: ApeXProtocol
$ws.Range('B44').Value = 'ThetaToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D44').Value = '''3.08'
$ws.Range('E44').Value = '  -2.64%  '
# Row 45: ThetaToken
$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').Value = '''3.61'
$ws.Range('E45').Value = '  +7.69%  '
# Row 46: VeChain
$ws.Range('D46').Value = '''0.0457'
$ws.Range('E46').Value = '  -0.24%  '
# Row 47: Maker
$ws.Range('D47').Value = '3.211.98'
$ws.Range('E47').Value = '  -3.20%  '
# Row 48: THORChain
$ws.Range('D48').Value = '''9.45'
$ws.Range('E48').Value = '  -2.30%  '
# Row 49: Stellar
$ws.Range('D49').Value = '''0.136'
$ws.Range('E49').Value = '  -1.89%  '
# Row 50: FirstDigitalUSD
$ws.Range('D50').Value = '''0.999'
$ws.Range('E50').Value = '  -0.05%  '
# Row 51: OceanProtocol
$ws.Range('D51').Value = '''1.45'
$ws.Range('E51').Value = '  +20.58%  '
